# The workbook has 529 data rows (rows 1-529), each with a value "LL" in
# column H ("PusMstKs" policy/status code). The fix corrects the typo
# "LL" -> "LU" for every row in column H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 529

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Value2 -eq "LL") {
        $cell.Value = "LU"
    }
}
